$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.838.39'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '1.561.55'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'205.74"
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('E6').Value = '  -2.04%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'21.70"
$ws.Range('E8').Value = '  -2.09%  '
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.783.44'
$ws.Range('E12').Value = '  -0.69%  '
$ws.Range('D13').Value = '1.562.64'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '26.834.47'
$ws.Range('E16').Value = '  -2.08%  '
$ws.Range('D17').Value = "'61.08"
$ws.Range('E17').Value = '  -3.73%  '
$ws.Range('D18').Value = "'213.86"
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = "'7.34"
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').Value = '0.0₃0678'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').Value = "'153.78"
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = "'6.71"
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').Value = "'14.93"
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('D31').Value = "'1.10"
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').Value = '1.399.67'
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').Value = "'1.52"
$ws.Range('E35').Value = '  -2.25%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('D40').Value = "'0.811"
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  -0.41%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'5.30"
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = "'2.17"
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'1.76"
$ws.Range('E45').Value = '  -2.52%  '
$ws.Range('D46').Value = "'62.93"
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').Value = '1.696.71'
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('D48').Value = "'85.96"
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').Value = '0.0₇0985'
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('D50').Value = "'0.0502"
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('D51').Value = "'0.0943"
$ws.Range('E51').Value = '  -0.72%  '

Write-Host "Applied 83 cell updates"
